$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.954.81'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '2.481.09'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'316.30"
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').Value = "'105.03"
$ws.Range('E6').Value = '  -4.02%  '
$ws.Range('D7').Value = "'0.518"
$ws.Range('E7').Value = '  -3.15%  '
$ws.Range('D9').Value = "'0.536"
$ws.Range('E9').Value = '  -3.84%  '
$ws.Range('D10').Value = "'38.76"
$ws.Range('E10').Value = '  -5.45%  '
$ws.Range('D11').Value = "'20.08"
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D12').Value = "'0.0799"
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').Value = "'7.06"
$ws.Range('E14').Value = '  -3.06%  '
$ws.Range('D15').Value = '2.871.24'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '2.484.27'
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('D17').Value = "'0.825"
$ws.Range('E17').Value = '  -3.88%  '
$ws.Range('D18').Value = '47.884.37'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  +8.87%  '
$ws.Range('D20').Value = "'12.66"
$ws.Range('E20').Value = '  -4.29%  '
$ws.Range('D21').Value = "'6.53"
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('D22').Value = '0.0₃0927'
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('D23').Value = "'70.68"
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('D24').Value = "'270.80"
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('D25').Value = "'2.50"
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('D26').Value = "'0.998"
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = "'25.61"
$ws.Range('E27').Value = '  -2.34%  '
$ws.Range('D28').Value = "'2.27"
$ws.Range('E28').Value = '  +6.79%  '
$ws.Range('D29').Value = "'9.66"
$ws.Range('E29').Value = '  -4.92%  '
$ws.Range('D30').Value = "'0.138"
$ws.Range('D31').Value = "'34.47"
$ws.Range('E31').Value = '  -4.58%  '
$ws.Range('D32').Value = "'49.33"
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = "'18.91"
$ws.Range('E34').Value = '  -5.43%  '
$ws.Range('D35').Value = "'5.24"
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('D37').Value = "'1.92"
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('D38').Value = "'4.53"
$ws.Range('E38').Value = '  -4.54%  '
$ws.Range('D39').Value = "'2.85"
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('D40').Value = "'122.35"
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').Value = "'22.10"
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').Value = "'0.0302"
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '1.997.82'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').Value = "'3.14"
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').Value = "'1.91"
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('E49').Value = '  -3.05%  '
$ws.Range('D50').Value = "'5.14"
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('D51').Value = "'78.26"
$ws.Range('E51').Value = '  -1.58%  '
